# PrefArr.xlsx - refresh the preference-array sample data and extend it
# from 34 to 43 rows, then leave the view scrolled/selected on the new
# data (matching the author's last editing session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in column A (rows 1-34) ---
# Only the rows whose value actually changed are listed; the rest keep
# their original contents untouched.
$ws.Range("A4").Value  = 0
$ws.Range("A5").Value  = 2
$ws.Range("A9").Value  = 0
$ws.Range("A10").Value = 0
$ws.Range("A13").Value = 1
$ws.Range("A14").Value = 1
$ws.Range("A15").Value = 0
$ws.Range("A16").Value = 1
$ws.Range("A21").Value = 1

# --- Append new rows 35-43 with value 0, extending the used range ---
$ws.Range("A35").Value = 0
$ws.Range("A36").Value = 0
$ws.Range("A37").Value = 0
$ws.Range("A38").Value = 0
$ws.Range("A39").Value = 0
$ws.Range("A40").Value = 0
$ws.Range("A41").Value = 0
$ws.Range("A42").Value = 0
$ws.Range("A43").Value = 0

# --- Update the saved view/selection state ---
# Scroll back to the top of the sheet and leave the selection on C40,
# mirroring where the author finished editing.
$ws.Range("A1").Select()
$ws.Range("C40").Select()
